# Menu-Languages.docx (Albanian) - rename the "RPC Explorer" menu entry
# to "Insight Explorer" (the resource link text itself is unchanged,
# only the visible run text differs from the source diff).
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "RPC Explorer",   # FindText
    $true,            # MatchCase
    $false,           # MatchWholeWord
    $false,           # MatchWildcards
    $false,           # MatchSoundsLike
    $false,           # MatchAllWordForms
    $true,            # Forward
    1,                # Wrap (wdFindContinue)
    $false,           # Format
    "Insight Explorer", # ReplaceWith
    2                 # Replace (wdReplaceAll)
)
